$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BuildingAgents reformatted to NetConnection Agents:
# column B ("agenttype") values of "building" become "netConnection"
# for the building-agent rows (rows 7 through 33).
for ($r = 7; $r -le 33; $r++) {
    $ws.Cells.Item($r, 2).Value = "netConnection"
}

# Reflect the resulting selection/active cell (B33) as in the saved workbook
$ws.Range("B33").Select()
